$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -6.565099999999998
$ws.Range("E4").Value = 13.3902

$ws.Range("E5").Value = 13.3412

$ws.Range("D6").Value = -7.789899999999998
$ws.Range("E6").Value = 13.8844

$ws.Range("D7").Value = -7.449799999999991

$ws.Range("D8").Value = -7.516199999999993
$ws.Range("E8").Value = 14.36099999999999

$ws.Range("D16").Value = -7.660599999999996
$ws.Range("E16").Value = 13.9821

$ws.Range("D20").Value = -7.571799999999999

$ws.Range("D21").Value = -7.590799999999996

$ws.Range("E22").Value = 13.37989999999999
